# Updates cryptos.xlsx cell values to match the "Updated symbol list" refresh
# (GitHub Actions run on Sat Dec 24 23:10:25 UTC 2022).
#
# Price/volume/hour columns hold numbers formatted as plain text in the source
# data (e.g. "244.48", "4.000", "23"), so for any value that looks numeric we
# force the cell to Text format before assigning it -- this prevents Excel from
# auto-converting the string into a real number (which would silently drop
# trailing zeros / introduce floating point noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether it must be kept as text
$updates = @(
    @{ Cell = "D2"; Value = "244.48"; ForceText = $true }
    @{ Cell = "G2"; Value = "23"; ForceText = $true }
    @{ Cell = "D3"; Value = "21.86"; ForceText = $true }
    @{ Cell = "G3"; Value = "23"; ForceText = $true }
    @{ Cell = "D4"; Value = "5.406"; ForceText = $true }
    @{ Cell = "G4"; Value = "23"; ForceText = $true }
    @{ Cell = "D5"; Value = "0.06035"; ForceText = $true }
    @{ Cell = "G5"; Value = "23"; ForceText = $true }
    @{ Cell = "D6"; Value = "3.398"; ForceText = $true }
    @{ Cell = "G6"; Value = "23"; ForceText = $true }
    @{ Cell = "D7"; Value = "0.8139"; ForceText = $true }
    @{ Cell = "G7"; Value = "23"; ForceText = $true }
    @{ Cell = "D8"; Value = "0.9240"; ForceText = $true }
    @{ Cell = "G8"; Value = "23"; ForceText = $true }
    @{ Cell = "D9"; Value = "0.1441"; ForceText = $true }
    @{ Cell = "G9"; Value = "23"; ForceText = $true }
    @{ Cell = "D10"; Value = "0.07472"; ForceText = $true }
    @{ Cell = "G10"; Value = "23"; ForceText = $true }
    @{ Cell = "D11"; Value = "0.03391"; ForceText = $true }
    @{ Cell = "G11"; Value = "23"; ForceText = $true }
    @{ Cell = "D12"; Value = "0.03069"; ForceText = $true }
    @{ Cell = "G12"; Value = "23"; ForceText = $true }
    @{ Cell = "D13"; Value = "0.09415"; ForceText = $true }
    @{ Cell = "G13"; Value = "23"; ForceText = $true }
    @{ Cell = "D14"; Value = "4.000"; ForceText = $true }
    @{ Cell = "G14"; Value = "23"; ForceText = $true }
    @{ Cell = "D15"; Value = "0.001588"; ForceText = $true }
    @{ Cell = "G15"; Value = "23"; ForceText = $true }
    @{ Cell = "D16"; Value = "0.04811"; ForceText = $true }
    @{ Cell = "G16"; Value = "23"; ForceText = $true }
    @{ Cell = "D17"; Value = "0.0005944"; ForceText = $true }
    @{ Cell = "G17"; Value = "23"; ForceText = $true }
    @{ Cell = "D18"; Value = "0.005639"; ForceText = $true }
    @{ Cell = "G18"; Value = "23"; ForceText = $true }
    @{ Cell = "D19"; Value = "0.004153"; ForceText = $true }
    @{ Cell = "G19"; Value = "23"; ForceText = $true }
    @{ Cell = "D20"; Value = "0.0009911"; ForceText = $true }
    @{ Cell = "G20"; Value = "23"; ForceText = $true }
    @{ Cell = "B21"; Value = "NitroEx"; ForceText = $false }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"; ForceText = $false }
    @{ Cell = "D21"; Value = "0.00008806"; ForceText = $true }
    @{ Cell = "E21"; Value = "20NitroExNTX"; ForceText = $false }
    @{ Cell = "G21"; Value = "23"; ForceText = $true }
    @{ Cell = "B22"; Value = "LEO"; ForceText = $false }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; ForceText = $false }
    @{ Cell = "D22"; Value = "3.652"; ForceText = $true }
    @{ Cell = "E22"; Value = "21LEOLEO"; ForceText = $false }
    @{ Cell = "G22"; Value = "23"; ForceText = $true }
    @{ Cell = "B23"; Value = "KuCoinToken"; ForceText = $false }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; ForceText = $false }
    @{ Cell = "D23"; Value = "6.429"; ForceText = $true }
    @{ Cell = "E23"; Value = "22KuCoinTokenKCS"; ForceText = $false }
    @{ Cell = "G23"; Value = "23"; ForceText = $true }
    @{ Cell = "B24"; Value = "BTSEToken"; ForceText = $false }
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; ForceText = $false }
    @{ Cell = "D24"; Value = "2.144"; ForceText = $true }
    @{ Cell = "E24"; Value = "23BTSETokenBTSE"; ForceText = $false }
    @{ Cell = "G24"; Value = "23"; ForceText = $true }
    @{ Cell = "B25"; Value = "BitpandaEcosystemToken"; ForceText = $false }
    @{ Cell = "C25"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; ForceText = $false }
    @{ Cell = "D25"; Value = "0.3243"; ForceText = $true }
    @{ Cell = "E25"; Value = "24BitpandaEcosystemTokenBEST"; ForceText = $false }
    @{ Cell = "G25"; Value = "23"; ForceText = $true }
    @{ Cell = "B26"; Value = "ProBitToken"; ForceText = $false }
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; ForceText = $false }
    @{ Cell = "D26"; Value = "0.1323"; ForceText = $true }
    @{ Cell = "E26"; Value = "25ProBitTokenPROB"; ForceText = $false }
    @{ Cell = "G26"; Value = "23"; ForceText = $true }
    @{ Cell = "D27"; Value = "0.0002902"; ForceText = $true }
    @{ Cell = "E27"; Value = "26UpBotsUBXT"; ForceText = $false }
    @{ Cell = "G27"; Value = "23"; ForceText = $true }
    @{ Cell = "G28"; Value = "23"; ForceText = $true }
    @{ Cell = "G29"; Value = "23"; ForceText = $true }
    @{ Cell = "G30"; Value = "23"; ForceText = $true }
    @{ Cell = "G31"; Value = "23"; ForceText = $true }
    @{ Cell = "G32"; Value = "23"; ForceText = $true }
    @{ Cell = "G33"; Value = "23"; ForceText = $true }
    @{ Cell = "G34"; Value = "23"; ForceText = $true }
    @{ Cell = "G35"; Value = "23"; ForceText = $true }
    @{ Cell = "G36"; Value = "23"; ForceText = $true }
    @{ Cell = "G37"; Value = "23"; ForceText = $true }
    @{ Cell = "G38"; Value = "23"; ForceText = $true }
    @{ Cell = "G39"; Value = "23"; ForceText = $true }
    @{ Cell = "D40"; Value = "0.03992"; ForceText = $true }
    @{ Cell = "G40"; Value = "23"; ForceText = $true }
    @{ Cell = "B41"; Value = "KickToken"; ForceText = $false }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; ForceText = $false }
    @{ Cell = "D41"; Value = "0.006434"; ForceText = $true }
    @{ Cell = "E41"; Value = "40KickTokenKICKBestin24h"; ForceText = $false }
    @{ Cell = "G41"; Value = "23"; ForceText = $true }
    @{ Cell = "B42"; Value = "BKEXToken"; ForceText = $false }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; ForceText = $false }
    @{ Cell = "D42"; Value = "0.1074"; ForceText = $true }
    @{ Cell = "E42"; Value = "41BKEXTokenBKK"; ForceText = $false }
    @{ Cell = "G42"; Value = "23"; ForceText = $true }
    @{ Cell = "B43"; Value = "CEJI"; ForceText = $false }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; ForceText = $false }
    @{ Cell = "D43"; Value = "0.002902"; ForceText = $true }
    @{ Cell = "E43"; Value = "42CEJICEJI"; ForceText = $false }
    @{ Cell = "G43"; Value = "23"; ForceText = $true }
    @{ Cell = "D44"; Value = "0.006391"; ForceText = $true }
    @{ Cell = "G44"; Value = "23"; ForceText = $true }
    @{ Cell = "D45"; Value = "0.00005250"; ForceText = $true }
    @{ Cell = "G45"; Value = "23"; ForceText = $true }
    @{ Cell = "D46"; Value = "0.00000000751"; ForceText = $true }
    @{ Cell = "G46"; Value = "23"; ForceText = $true }
    @{ Cell = "D47"; Value = "1.101"; ForceText = $true }
    @{ Cell = "G47"; Value = "23"; ForceText = $true }
    @{ Cell = "D48"; Value = "0.002320"; ForceText = $true }
    @{ Cell = "G48"; Value = "23"; ForceText = $true }
    @{ Cell = "G49"; Value = "23"; ForceText = $true }
    @{ Cell = "D50"; Value = "0.01011"; ForceText = $true }
    @{ Cell = "G50"; Value = "23"; ForceText = $true }
    @{ Cell = "G51"; Value = "23"; ForceText = $true }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $update.Value
}

